# agrego fechas de proyecto final BD
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "Plan de estudio"

# --- Column M width: 33.85546875 -> 40 (character units) ---
$ws.Range("M1").EntireColumn.ColumnWidth = 39.166666666666664

# --- Update existing "Dias para la entrega" counts (col L) ---
$ws.Range("L3").Value2 = 0
$ws.Range("L4").Value2 = 2
$ws.Range("L5").Value2 = 5
$ws.Range("L6").Value2 = 5

# --- Rows 7-9: new "PROYECTO FINAL" deliveries ---
# Copy the format from row 3's "Pendiente" cell (style 12) onto O7:O9 first,
# so the new values pick up the same fill/border/alignment used by the
# other "Pendiente" rows.
$ws.Range("O3").Copy() | Out-Null
$ws.Range("O7:O9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Column J (materia) - reuses the existing "Bases de datos" shared string.
$ws.Range("J7").Value2 = "Bases de datos"
$ws.Range("J8").Value2 = "Bases de datos"
$ws.Range("J9").Value2 = "Bases de datos"

# Column K (fecha de entrega) - new shared strings appended in this order
# (83, 84, 85) to mirror the original commit's sharedStrings.xml ordering.
$ws.Range("K7").Value2 = "Viernes 5 de junio"
$ws.Range("K8").Value2 = "Viernes 10 de julio"
$ws.Range("K9").Value2 = "viernes 7 de agosto"

# Column L (dias para la entrega)
$ws.Range("L7").Value2 = 10

# Column M (descripcion de tarea) - new shared strings appended next
# (86, 87, 88).
$ws.Range("M7").Value2 = "Primer entrega PROYECTO FINAL"
$ws.Range("M8").Value2 = "Segunda entrega PROYECTO FINAL"
$ws.Range("M9").Value2 = "Tercer entrega PROYECTO FINAL"

# Column N (nivel de importancia)
$ws.Range("N7").Value2 = 5
$ws.Range("N8").Value2 = 5
$ws.Range("N9").Value2 = 5

# Column O (estado) - reuses the existing "Pendiente" shared string.
$ws.Range("O7").Value2 = "Pendiente"
$ws.Range("O8").Value2 = "Pendiente"
$ws.Range("O9").Value2 = "Pendiente"

$excel.CutCopyMode = $false

# --- Selection moves to M10 ---
$ws.Activate() | Out-Null
$ws.Range("M10").Select() | Out-Null
